$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '94.009.76'
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').Value = '3.060.25'
$ws.Range('E3').Value = '  -2.08%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.14'
$ws.Range('E5').Value = '  -4.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '604.83'
$ws.Range('E6').Value = '  -1.88%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.08'
$ws.Range('E7').Value = '  -2.06%  '
$ws.Range('E8').Value = '  -7.89%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('E10').Value = '  +7.03%  '
$ws.Range('D11').Value = '3.057.77'
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('E12').Value = '  -4.00%  '
$ws.Range('D13').Value = '93.665.26'
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000237'
$ws.Range('E14').Value = '  -6.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '33.28'
$ws.Range('E15').Value = '  -4.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.24'
$ws.Range('E16').Value = '  -4.47%  '
$ws.Range('D17').Value = '3.629.51'
$ws.Range('E17').Value = '  -2.63%  '
$ws.Range('D18').Value = '3.046.18'
$ws.Range('E18').Value = '  -3.26%  '
$ws.Range('E19').Value = '  -8.56%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.19'
$ws.Range('E20').Value = '  -4.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.61'
$ws.Range('E21').Value = '  -3.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '434.96'
$ws.Range('E22').Value = '  -3.23%  '
$ws.Range('E23').Value = '  -7.41%  '
$ws.Range('E24').Value = '  -10.23%  '
$ws.Range('E25').Value = '  +5.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.43'
$ws.Range('E26').Value = '  -7.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '84.19'
$ws.Range('E27').Value = '  -3.82%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.54'
$ws.Range('E28').Value = '  -2.44%  '
$ws.Range('D29').Value = '3.222.70'
$ws.Range('E29').Value = '  -2.29%  '
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.242'
$ws.Range('E31').Value = '  +5.02%  '
$ws.Range('B32').Value = 'Cronos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.175'
$ws.Range('E32').Value = '  +2.87%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.121'
$ws.Range('E33').Value = '  -12.16%  '
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.02'
$ws.Range('E34').Value = '  +2.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.91'
$ws.Range('E35').Value = '  -4.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.56'
$ws.Range('E36').Value = '  -6.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.154'
$ws.Range('E37').Value = '  -3.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.08'
$ws.Range('E38').Value = '  -4.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.87'
$ws.Range('E39').Value = '  -2.64%  '
$ws.Range('E40').Value = '  +3.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.76'
$ws.Range('E41').Value = '  -2.69%  '
$ws.Range('E42').Value = '  -1.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '457.57'
$ws.Range('E43').Value = '  -7.87%  '
$ws.Range('E44').Value = '  -5.86%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.07'
$ws.Range('E46').Value = '  -10.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '160.09'
$ws.Range('E47').Value = '  -2.24%  '
$ws.Range('E48').Value = '  -7.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.657'
$ws.Range('E49').Value = '  -5.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '43.59'
$ws.Range('E50').Value = '  -1.16%  '
$ws.Range('E51').Value = '  +0.01%  '
